$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---
$ws.Range("D2").Value = 260
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = -1
$ws.Range("G2").Value = -7
$ws.Range("H2").Value = -7
$ws.Range("I2").Value = -7
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 304
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 259
$ws.Range("N2").Value = 259
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 142
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = -14
$ws.Range("S2").Value = 23
$ws.Range("T2").Value = 16
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 20
$ws.Range("W2").Value = -0.53
$ws.Range("X2").Value = -2.83
$ws.Range("Y2").Value = -3.62
$ws.Range("Z2").Value = -2.5
$ws.Range("AA2").Value = 17.53
$ws.Range("AB2").Value = 77.42
$ws.Range("AC2").Value = -31
$ws.Range("AD2").Value = -66.01000000000001
$ws.Range("AE2").Value = 977
$ws.Range("AF2").Value = 2.06
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 26460000

# --- Row 3 (2015/12) ---
$ws.Range("D3").Value = 232
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 411
$ws.Range("L3").Value = 129
$ws.Range("M3").Value = 281
$ws.Range("N3").Value = 281
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 142
$ws.Range("Q3").Value = 31
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 78
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 28
$ws.Range("V3").Value = 95
$ws.Range("W3").Value = 1.81
$ws.Range("X3").Value = 8.970000000000001
$ws.Range("Y3").Value = 7.76
$ws.Range("Z3").Value = 5.82
$ws.Range("AA3").Value = 45.97
$ws.Range("AB3").Value = 91.86
$ws.Range("AC3").Value = 79
$ws.Range("AD3").Value = 87.42
$ws.Range("AE3").Value = 1063
$ws.Range("AF3").Value = 6.51
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 26460000

# --- Row 4 (2016/12) ---
$ws.Range("D4").Value = 207
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 11
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 411
$ws.Range("L4").Value = 56
$ws.Range("M4").Value = 355
$ws.Range("N4").Value = 355
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 156
$ws.Range("Q4").Value = 32
$ws.Range("R4").Value = -15
$ws.Range("S4").ClearContents()
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 28
$ws.Range("V4").Value = 35
$ws.Range("W4").Value = 10.22
$ws.Range("X4").Value = 5.26
$ws.Range("Y4").Value = 3.43
$ws.Range("Z4").Value = 2.65
$ws.Range("AA4").Value = 15.9
$ws.Range("AB4").Value = 122.63
$ws.Range("AC4").Value = 40
$ws.Range("AD4").Value = 102.56
$ws.Range("AE4").Value = 1210
$ws.Range("AF4").Value = 3.35
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 29317140

# --- Row 5 (2017/12) ---
$ws.Range("D5").Value = 187
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 13
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 404
$ws.Range("L5").Value = 21
$ws.Range("M5").Value = 383
$ws.Range("N5").Value = 383
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 159
$ws.Range("Q5").Value = 12
$ws.Range("R5").Value = -93
$ws.Range("S5").Value = -20
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 8
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 7.11
$ws.Range("X5").Value = 5.54
$ws.Range("Y5").Value = 2.81
$ws.Range("Z5").Value = 2.54
$ws.Range("AA5").Value = 5.47
$ws.Range("AB5").Value = 134.71
$ws.Range("AC5").Value = 35
$ws.Range("AD5").Value = 81.77
$ws.Range("AE5").Value = 1278
$ws.Range("AF5").Value = 2.22
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 29976480

# --- Row 6 (2018/12) ---
$ws.Range("D6").Value = 187
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 38
$ws.Range("I6").Value = 38
$ws.Range("K6").Value = 546
$ws.Range("L6").Value = 125
$ws.Range("M6").Value = 421
$ws.Range("N6").Value = 421
$ws.Range("P6").Value = 159
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = -41
$ws.Range("S6").Value = 98
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 77
$ws.Range("W6").Value = 3.46
$ws.Range("X6").Value = 20.29
$ws.Range("Y6").Value = 9.42
$ws.Range("Z6").Value = 7.98
$ws.Range("AA6").Value = 29.66
$ws.Range("AB6").Value = 164.26
$ws.Range("AC6").Value = 126
$ws.Range("AD6").Value = 18.36
$ws.Range("AE6").Value = 1405
$ws.Range("AF6").Value = 1.65
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 29976480

# --- Rows 7, 8, 9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# These projected-year rows lose all of their estimate figures (columns D
# through AI), leaving only the row index / ticker / period label columns.
$ws.Range("D7:AI9").ClearContents()
